$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Style bookkeeping -----------------------------------------------------
# E2 is currently the sole user of its own cell style (center/center, thin
# border, no wrap). Restyle it in place to bold text on a yellow fill (and
# add wrap) so that this cell's existing style slot is rewritten rather than
# a brand-new one being appended.
$ws.Range("E2").Font.Bold = $true
$ws.Range("E2").Interior.Color = 65535
$ws.Range("E2").WrapText = $true

# Apply the very same look (bold + yellow fill + centered/wrapped) to the
# last two header cells; this matches the style slot just rewritten above,
# so it gets reused instead of creating a new one.
$ws.Range("I1:J1").Font.Bold = $true
$ws.Range("I1:J1").Interior.Color = 65535

# --- Clear the sample data row ---------------------------------------------
$ws.Range("A2:J2").ClearContents()

# Give the bulk of the (now empty) data row a plain bordered look (no special
# alignment) - this covers the former distinctive E2 cell as well.
$ws.Range("C2:J2").ClearFormats()
$ws.Range("C2:J2").Borders.ColorIndex = 1
$ws.Range("C2:J2").Borders.LineStyle = 1
$ws.Range("C2:J2").Borders.Weight = 2

# --- Update the remembered selection ---------------------------------------
$ws.Range("D6").Select() | Out-Null
